$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4497.5
$ws.Range("I69").Value = 3996.6667
$ws.Range("K69").Value = 11990.0001
$ws.Range("M69").Value = -11116.0001
$ws.Range("H72").Value = 4497.5
$ws.Range("I72").Value = 3996.6667
$ws.Range("K72").Value = 35970.0003
$ws.Range("M72").Value = -31602.0003
$ws.Range("H111").Value = 5357
$ws.Range("I111").Value = 5416.6665
$ws.Range("K111").Value = 16249.9995
$ws.Range("M111").Value = -13182.9995
$ws.Range("H129").Value = 16345.857
$ws.Range("I129").Value = 1799.6666
$ws.Range("J129").Value = 27255.5
$ws.Range("K129").Value = 5398.9998
$ws.Range("L129").Value = 81766.5
$ws.Range("M129").Value = -398.9997999999996
$ws.Range("N129").Value = -91766.5
$ws.Range("H132").Value = 2410.6326
$ws.Range("I132").Value = 2481.4348
$ws.Range("J132").Value = 1325
$ws.Range("K132").Value = 7444.3044
$ws.Range("L132").Value = 3975
$ws.Range("M132").Value = -4914.3044
$ws.Range("N132").Value = -9035
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140
$ws.Range("H138").Value = 18871950
$ws.Range("I138").Value = 925.73914
$ws.Range("K138").Value = 2777.21742
$ws.Range("M138").Value = 2362.78258
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H37").Value = 31666.666
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40546
$ws.Range("H61").Value = 3302.158
$ws.Range("I61").Value = 3221.5881
$ws.Range("K61").Value = 3221.5881
$ws.Range("M61").Value = -3009.5881
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("H97").Value = 1645.8667
$ws.Range("I97").Value = 1639.909
$ws.Range("K97").Value = 1639.909
$ws.Range("M97").Value = -1143.909
$ws.Range("H122").Value = 1241.4348
$ws.Range("I122").Value = 1069.4103
$ws.Range("J122").Value = 2199.8572
$ws.Range("K122").Value = 3208.2309
$ws.Range("L122").Value = 6599.571599999999
$ws.Range("M122").Value = -758.2309
$ws.Range("N122").Value = -11499.5716
$ws.Range("H136").Value = 3302.158
$ws.Range("I136").Value = 3221.5881
$ws.Range("K136").Value = 9664.764299999999
$ws.Range("M136").Value = -7114.764299999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1383.6428
$ws.Range("I94").Value = 1194.2
$ws.Range("K94").Value = 1194.2
$ws.Range("M94").Value = -743.2
$ws.Range("H95").Value = 24839.143
$ws.Range("J95").Value = 24839.143
$ws.Range("L95").Value = 24839.143
$ws.Range("N95").Value = -30331.143
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4834.25
$ws.Range("J99").Value = 3746.75
$ws.Range("L99").Value = 3746.75
$ws.Range("N99").Value = -6742.75
$ws.Range("H126").Value = 4834.25
$ws.Range("J126").Value = 3746.75
$ws.Range("L126").Value = 11240.25
$ws.Range("N126").Value = -16180.25
$ws.Range("H130").Value = 57175
$ws.Range("J130").Value = 57175
$ws.Range("L130").Value = 57175
$ws.Range("N130").Value = -67215
$ws.Range("H134").Value = 1402.7727
$ws.Range("I134").Value = 1207.9474
$ws.Range("K134").Value = 3623.8422
$ws.Range("M134").Value = -1088.8422
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2004
$ws.Range("I62").Value = 2004
$ws.Range("K62").Value = 6012
$ws.Range("M62").Value = -5326
$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 300000
$ws.Range("N64").Value = -300540
$ws.Range("H65").Value = 2004
$ws.Range("I65").Value = 2004
$ws.Range("K65").Value = 18036
$ws.Range("M65").Value = -14604
$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 300000
$ws.Range("N67").Value = -301872
$ws.Range("H69").Value = 1449.5
$ws.Range("J69").Value = 1399
$ws.Range("L69").Value = 4197
$ws.Range("N69").Value = -5819
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H72").Value = 1449.5
$ws.Range("J72").Value = 1399
$ws.Range("L72").Value = 12591
$ws.Range("N72").Value = -20703
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H74").Value = 13663.333
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -47122
$ws.Range("H77").Value = 13663.333
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -145608
$ws.Range("H81").Value = 2501853
$ws.Range("I81").Value = 3334137.2
$ws.Range("K81").Value = 10002411.6
$ws.Range("M81").Value = -10001288.6
$ws.Range("H84").Value = 2501853
$ws.Range("I84").Value = 3334137.2
$ws.Range("K84").Value = 30007234.8
$ws.Range("M84").Value = -30001618.8
$ws.Range("H113").Value = 1096.75
$ws.Range("I113").Value = 403.85715
$ws.Range("K113").Value = 1211.57145
$ws.Range("M113").Value = 958.4285500000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1360.5483
$ws.Range("I122").Value = 1005.55554
$ws.Range("K122").Value = 3016.66662
$ws.Range("M122").Value = -566.66662
$ws.Range("H126").Value = 3689.5
$ws.Range("I126").Value = 2757.4285
$ws.Range("K126").Value = 8272.2855
$ws.Range("M126").Value = -5802.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1072.85
$ws.Range("I16").Value = 968.05884
$ws.Range("K16").Value = 968.05884
$ws.Range("M16").Value = -798.05884
$ws.Range("H132").Value = 1471.3125
$ws.Range("J132").Value = 1255
$ws.Range("L132").Value = 3765
$ws.Range("N132").Value = -8825
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 1699666.5
$ws.Range("J124").Value = 1699666.5
$ws.Range("L124").Value = 1699666.5
$ws.Range("N124").Value = -1709486.5
$ws.Range("H140").Value = 57664.5
$ws.Range("J140").Value = 57664.5
$ws.Range("L140").Value = 57664.5
$ws.Range("N140").Value = -68024.5
